$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 3.9
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("S4").Value = 2.88
$ws.Range("T4").Value = 1.4
$ws.Range("U4").Value = 4.7
$ws.Range("V4").Value = 1.19
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 1.13
$ws.Range("AA4").Value = 2.38
$ws.Range("AB4").Value = 1.53
$ws.Range("AC4").Value = 5
$ws.Range("AE4").Value = 10
$ws.Range("AK4").Value = 23
$ws.Range("AL4").Value = 101
$ws.Range("AO4").Value = 19

# Row 5 updates
$ws.Range("G5").Value = 2.4
$ws.Range("J5").Value = 3.25
$ws.Range("K5").Value = 1.91
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 2.38
$ws.Range("S5").Value = 2.7
$ws.Range("T5").Value = 1.44
$ws.Range("U5").Value = 4.5
$ws.Range("V5").Value = 1.21
$ws.Range("W5").Value = 5.5
$ws.Range("X5").Value = 1.14
$ws.Range("AE5").Value = 10
$ws.Range("AG5").Value = 23
$ws.Range("AK5").Value = 19
$ws.Range("AN5").Value = 7.5
$ws.Range("AS5").Value = 41
